$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update importance values in column B (rows 2-10)
$ws.Range("B2").Value = 0.4068687770768175
$ws.Range("B3").Value = 0.1849681479539378
$ws.Range("B4").Value = 0.09077552765695136
$ws.Range("B5").Value = 0.0865054058318841
$ws.Range("B6").Value = 0.05884279686526506
$ws.Range("B7").Value = 0.05328526475793372
$ws.Range("B8").Value = 0.05144639809302326
$ws.Range("B9").Value = 0.04760883173027003
$ws.Range("B10").Value = 0.01969885003391711

# Update reordered feature labels in column A (rows 5-9)
$ws.Range("A5").Value = "VIX_short"
$ws.Range("A6").Value = "VIX_long"
$ws.Range("A7").Value = "VIX"
$ws.Range("A8").Value = "close_long"
$ws.Range("A9").Value = "close_short"
